# The frequency column (column A, rows 2-17) on both data sheets was
# mistakenly written in rad/s instead of Hz. Convert it by dividing every
# value by 2*pi (omega = 2*pi*f  =>  f = omega / (2*pi)).
$twoPi = 6.283185307179586

$wb = $excel.ActiveWorkbook

$sheet180 = $wb.Worksheets.Item("180")
$sheet160 = $wb.Worksheets.Item("160")

# --- Convert frequency values on sheet "180" ---
for ($r = 2; $r -le 17; $r++) {
    $cell = $sheet180.Cells.Item($r, 1)
    $old = $cell.Value()
    $cell.Value = $old / $twoPi
}

# --- Convert frequency values on sheet "160" ---
for ($r = 2; $r -le 17; $r++) {
    $cell = $sheet160.Cells.Item($r, 1)
    $old = $cell.Value()
    $cell.Value = $old / $twoPi
}

# --- Update the selection / active sheet state to match the saved view ---
# Sheet "160" ends up with A2:A17 selected (not the active/visible tab).
$sheet160.Select()
$sheet160.Range("A2:A17").Select()

# Sheet "180" ends up active (tabSelected) with A2:A17 selected.
$sheet180.Select()
$sheet180.Range("A2:A17").Select()
